# customFunction.xlsx: fix the "Formula Expression" label text so each row's
# displayed formula string matches the formula actually used in column C on
# that row (labels had drifted to the wrong row's arguments), and flip which
# sheet/cell is active/selected when the workbook is next opened.

$wb = $excel.ActiveWorkbook

$wsCustom = $wb.Worksheets.Item("custom-formula")
$wsChain  = $wb.Worksheets.Item("CHAIN")

# --- Correct the mislabeled "Formula Expression" text cells -----------------
$wsCustom.Range("B5").Value  = " =MYEXCHANGE(D5, E5)"
$wsCustom.Range("B7").Value  = " =MYSUBTOTAL(D7:G7)"
$wsCustom.Range("B9").Value  = " =MYSUBTOTAL(D9:G9)"
$wsCustom.Range("B11").Value = " =MYSUBTOTAL(D11:G11)"

$wsChain.Range("B4").Value = " =CHAIN(D4:G4)"
$wsChain.Range("B6").Value = " =CHAIN(D6:G6)"
$wsChain.Range("B8").Value = " =CHAIN(D8,E8,F8,G8)"

# --- Update selection / active sheet ----------------------------------------
# custom-formula keeps its own remembered selection (B13) but is no longer
# the tab shown on open; CHAIN becomes the active tab with B8 selected.
# (Selecting a range also makes its sheet the active one, so set the
# non-active sheet's selection first, then activate+select CHAIN last so it
# ends up as the visible/active tab.)
[void]$wsCustom.Range("B13").Select()

[void]$wsChain.Activate()
[void]$wsChain.Range("B8").Select()
